# crossbar - updated positions of the leds
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crossBar")

# --- Move LED_x / GPIO marker pairs (columns D/E) to their new pin rows ---
# LED_0: row 11 -> row 8
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("D8").Value = "LED_0"
$ws.Range("E8").Value = "GPIO"

# LED_1: stays logically but the source row (12) previously held nothing;
# LED_1 moves from row 14 -> row 12, and row 14 becomes LED_2
$ws.Range("D12").Value = "LED_1"
$ws.Range("E12").Value = "GPIO"

# LED_1 -> LED_2 in place at row 14 (value-only change)
$ws.Range("D14").Value = "LED_2"

# LED_3: row 23 -> row 18
$ws.Range("D23").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("D18").Value = "LED_3"
$ws.Range("E18").Value = "GPIO"

# LED_4: row 25 -> row 20
$ws.Range("D25").ClearContents()
$ws.Range("E25").ClearContents()
$ws.Range("D20").Value = "LED_4"
$ws.Range("E20").Value = "GPIO"

# old LED_2 at row 21 removed (superseded by LED_2 now at row 14)
$ws.Range("D21").ClearContents()
$ws.Range("E21").ClearContents()

# LED_5: row 27 -> row 22
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("D22").Value = "LED_5"
$ws.Range("E22").Value = "GPIO"

# LED_6: row 29 -> row 26
$ws.Range("D29").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("D26").Value = "LED_6"
$ws.Range("E26").Value = "GPIO"

# LED_7: row 31 -> row 28
$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()
$ws.Range("D28").Value = "LED_7"
$ws.Range("E28").Value = "GPIO"

# --- View state: frozen pane moved up one row, selection moved ---
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("E18").Select()
